$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 301, shifting the existing rows 301:323 down to 302:324
$ws.Rows("301:301").Insert()

# Populate the newly inserted row 301 with the new data record
$ws.Range("A301").Value = 4
$ws.Range("B301").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C301").Value = "Los Lagos"
$ws.Range("D301").Value = 44746
$ws.Range("E301").Value = 10
$ws.Range("F301").Value = 100112045
$ws.Range("G301").Value = "Zapallo"
$ws.Range("H301").Value = "Paine"
$ws.Range("I301").Value = "1a (guarda)"
$ws.Range("J301").Value = 500
$ws.Range("K301").Value = 500
$ws.Range("L301").Value = 500
$ws.Range("M301").Value = 500
$ws.Range("N301").Value = '$/kilo (volumen en unidades)'
$ws.Range("O301").Value = "Región de O'Higgins"
$ws.Range("P301").Value = 500
$ws.Range("Q301").Value = 1
$ws.Range("R301").Value = "Hortaliza"
